# Fruta / hortaliza, semanal
# Updates columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg)
# for rows 2-11 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row, taken directly from the target workbook state.
# Columns: D, L, M, N, O, P, S
$updates = @{
    2  = @{ D = 44511; L = "Primera"; M = 120; N = 28000; O = 28000; P = 28000; S = 2800 }
    3  = @{ D = 44432; L = "Primera"; M = 20;  N = 20000; O = 20000; P = 20000; S = 2000 }
    4  = @{ D = 44435; L = "Primera"; M = 40;  N = 20000; O = 20000; P = 20000; S = 2000 }
    5  = @{ D = 44503; L = "Primera"; M = 60;  N = 30000; O = 30000; P = 30000; S = 3000 }
    6  = @{ D = 44503; L = "Segunda"; M = 50;  N = 25000; O = 25000; P = 25000; S = 2500 }
    7  = @{ D = 44473; L = "Primera"; M = 180; N = 20000; O = 20000; P = 20000; S = 2000 }
    8  = @{ D = 44517; L = "Especial"; M = 100; N = 27000; O = 27000; P = 27000; S = 2700 }
    9  = @{ D = 44517; L = "Primera"; M = 30;  N = 25000; O = 25000; P = 25000; S = 2500 }
    10 = @{ D = 44434; L = "Primera"; M = 20;  N = 20000; O = 20000; P = 20000; S = 2000 }
    11 = @{ D = 44466; L = "Primera"; M = 60;  N = 20000; O = 20000; P = 20000; S = 2000 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
